# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to match the newly generated output, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> newValue } for column F
$sheetChanges = @{
    "展览" = @{
        3  = 514
        4  = 1506
        6  = 138
        8  = 140
        9  = 730
        10 = 1041
        11 = 61
        12 = 322
        13 = 49
        14 = 6351
        15 = 2
        18 = 147
        20 = 15222
        21 = 1508
        22 = 276
        23 = 136
        25 = 11013
        26 = 743
        27 = 4305
        28 = 228
        29 = 372
        30 = 15
    }
    "全部类型" = @{
        3  = 514
        4  = 1506
        6  = 138
        9  = 140
        10 = 730
        12 = 1041
        13 = 61
        14 = 322
        15 = 49
        17 = 6351
        18 = 2
        21 = 147
        23 = 15222
        24 = 1508
        25 = 276
        26 = 136
        28 = 11013
        29 = 743
        30 = 4305
        31 = 228
        32 = 372
        33 = 15
    }
}

foreach ($sheetName in $sheetChanges.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetChanges[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item([int]$row, 6).Value = $rows[$row]
    }
}
